$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Historical GC")

$ws.Range("C2").Value = "'8.2662"
$ws.Range("C3").Value = "'3.5367"
$ws.Range("C4").Value = "'11.3173"
$ws.Range("C5").Value = "'6.6169"
$ws.Range("C6").Value = "'4.8026"
$ws.Range("C7").Value = "'0.9605"
$ws.Range("C8").Value = "'8.1438"
$ws.Range("C9").Value = "'0.4309"
$ws.Range("C10").Value = "'4.8451"
$ws.Range("C11").Value = "'0.7673"
$ws.Range("C12").Value = "'5.2595"
$ws.Range("C13").Value = "'6.0340"
$ws.Range("C14").Value = "'0.9382"
$ws.Range("C15").Value = "'10.8191"
$ws.Range("C16").Value = "'1.7711"
$ws.Range("C17").Value = "'1.8826"
$ws.Range("C18").Value = "'5.3119"
$ws.Range("C19").Value = "'6.4580"
$ws.Range("C20").Value = "'1.7218"
$ws.Range("C21").Value = "'1.0898"
$ws.Range("C22").Value = "'5.3218"
$ws.Range("C23").Value = "'7.8135"
$ws.Range("C24").Value = "'4.8382"
$ws.Range("C25").Value = "'-0.8065"
$ws.Range("C26").Value = "'5.0152"
$ws.Range("C27").Value = "'8.2548"
$ws.Range("C28").Value = "'4.9719"
$ws.Range("C29").Value = "'5.5496"
$ws.Range("C30").Value = "'5.8369"
$ws.Range("C31").Value = "'-2.3341"
$ws.Range("C32").Value = "'5.1956"
$ws.Range("C33").Value = "'5.8853"
$ws.Range("C34").Value = "'-4.9243"
$ws.Range("C35").Value = "'11.2981"
$ws.Range("C36").Value = "'5.2365"
$ws.Range("C37").Value = "'4.7313"
$ws.Range("C38").Value = "'4.1272"
$ws.Range("C39").Value = "'4.5509"
$ws.Range("C40").Value = "'7.4017"
$ws.Range("C41").Value = "'5.4089"
$ws.Range("C42").Value = "'9.4576"
$ws.Range("C43").Value = "'6.0079"
$ws.Range("C44").Value = "'3.8367"
$ws.Range("C45").Value = "'5.8154"
$ws.Range("C46").Value = "'75.9336"
$ws.Range("C47").Value = "'6.4483"
$ws.Range("C48").Value = "'6.1999"
$ws.Range("C49").Value = "'10.9216"
$ws.Range("C50").Value = "'6.0997"
$ws.Range("C51").Value = "'20.1638"
$ws.Range("C52").Value = "'7.8909"
$ws.Range("C53").Value = "'6.9303"
$ws.Range("C54").Value = "'5.9085"
$ws.Range("C55").Value = "'5.1385"
$ws.Range("C56").Value = "'4.7512"
$ws.Range("C57").Value = "'-0.9414"
$ws.Range("C58").Value = "'5.3067"
$ws.Range("C59").Value = "'2.8141"
$ws.Range("C60").Value = "'4.8281"
$ws.Range("C61").Value = "'2.0201"
$ws.Range("C62").Value = "'5.2532"
$ws.Range("C63").Value = "'0.8278"
$ws.Range("C64").Value = "'6.8696"
$ws.Range("C65").Value = "'-18.8433"
$ws.Range("C66").Value = "'4.9600"
$ws.Range("C67").Value = "'2.8192"
$ws.Range("C68").Value = "'4.8994"
$ws.Range("C69").Value = "'3.6974"
$ws.Range("C70").Value = "'6.9790"
$ws.Range("C71").Value = "'-7.7848"
$ws.Range("C72").Value = "'5.9327"
$ws.Range("C73").Value = "'-49.3505"
$ws.Range("C74").Value = "'5.9085"
$ws.Range("C75").Value = "'6.2415"
$ws.Range("C76").Value = "'5.1417"
$ws.Range("C77").Value = "'-1.1973"
$ws.Range("C78").Value = "'5.1032"
$ws.Range("C79").Value = "'5.8866"
$ws.Range("C80").Value = "'1.5674"
$ws.Range("C81").Value = "'6.2041"
$ws.Range("C82").Value = "'8.1576"
$ws.Range("C83").Value = "'4.8026"
$ws.Range("C84").Value = "'0.9605"
$ws.Range("C85").Value = "'1.6797"
$ws.Range("C86").Value = "'5.4513"
$ws.Range("C87").Value = "'7.9952"
$ws.Range("C88").Value = "'0.3143"
$ws.Range("C89").Value = "'-0.2417"
$ws.Range("C90").Value = "'5.6647"
$ws.Range("C91").Value = "'-22.8655"
$ws.Range("C92").Value = "'5.0191"
$ws.Range("C93").Value = "'6.7117"
$ws.Range("C94").Value = "'5.4544"
$ws.Range("C95").Value = "'9.8971"
$ws.Range("C96").Value = "'5.3886"
$ws.Range("C97").Value = "'8.6206"
$ws.Range("C98").Value = "'5.2331"
$ws.Range("C99").Value = "'6.3303"
$ws.Range("C100").Value = "'-0.9588"
$ws.Range("C101").Value = "'6.0905"
$ws.Range("C102").Value = "'11.1399"
$ws.Range("C103").Value = "'6.0476"
$ws.Range("C104").Value = "'12.4009"
$ws.Range("C105").Value = "'5.1989"
$ws.Range("C106").Value = "'6.1437"
$ws.Range("C107").Value = "'6.6913"
$ws.Range("C108").Value = "'4.7785"
$ws.Range("C109").Value = "'5.0484"
$ws.Range("C110").Value = "'-0.8999"
$ws.Range("C111").Value = "'6.0904"
$ws.Range("C112").Value = "'-23.2214"
$ws.Range("C113").Value = "'5.9403"
$ws.Range("C114").Value = "'13.9284"
$ws.Range("C115").Value = "'6.4039"
$ws.Range("C116").Value = "'-38.1362"
$ws.Range("C117").Value = "'0.0187"
$ws.Range("C118").Value = "'6.1554"
$ws.Range("C119").Value = "'-32.2679"
$ws.Range("C120").Value = "'3.7460"
$ws.Range("C121").Value = "'5.7981"
$ws.Range("C122").Value = "'18.0214"
$ws.Range("C123").Value = "'7.6497"
$ws.Range("C124").Value = "'1.4492"
$ws.Range("C125").Value = "'6.0845"
$ws.Range("C126").Value = "'-28.2994"
$ws.Range("C127").Value = "'5.4008"
$ws.Range("C128").Value = "'6.1199"
$ws.Range("C129").Value = "'2.2302"
$ws.Range("C130").Value = "'4.7498"
$ws.Range("C131").Value = "'7.9329"
